$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 358, all hold the "Förändrad" (last changed) date
# serial value 45190 and need to be bumped to 45192.
$range = $ws.Range("C2:C358")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
